$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to be treated as text so values that look like
# plain numbers (e.g. "608.96") are not auto-converted to numeric cells -
# they must stay text like the rest of the column (e.g. "66.622.55").
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "66.605.88"
$ws.Range("E2").Value = "  +0.68%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.596.72"
$ws.Range("E3").Value = "  +0.96%  "

# Row 5 - BNB
$ws.Range("D5").Value = "608.96"
$ws.Range("E5").Value = "  +0.33%  "

# Row 6 - Solana
$ws.Range("D6").Value = "148.38"
$ws.Range("E6").Value = "  +2.38%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.08%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +0.31%  "

# Row 9 - was Toncoin, now Dogecoin
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "0.136"
$ws.Range("E9").Value = "  -0.13%  "

# Row 10 - was Dogecoin, now Toncoin
$ws.Range("B10").Value = "Toncoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D10").Value = "8.03"
$ws.Range("E10").Value = "  +0.32%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "4.206.40"
$ws.Range("E12").Value = "  +0.97%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  +0.72%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "29.77"
$ws.Range("E14").Value = "  -1.14%  "

# Row 15 - WrappedEther
$ws.Range("D15").Value = "3.586.94"
$ws.Range("E15").Value = "  +0.64%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "66.691.67"
$ws.Range("E16").Value = "  +0.66%  "

# Row 18 - Uniswap
$ws.Range("D18").Value = "11.50"
$ws.Range("E18").Value = "  +0.89%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  +2.00%  "

# Row 20 - Chainlink
$ws.Range("D20").Value = "15.10"
$ws.Range("E20").Value = "  +1.40%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "427.07"
$ws.Range("E21").Value = "  -0.92%  "

# Row 22 - Polygon
$ws.Range("E22").Value = "  +1.19%  "

# Row 23 - Litecoin
$ws.Range("E23").Value = "  +0.07%  "

# Row 24 - WrappedeETH
$ws.Range("D24").Value = "3.736.48"
$ws.Range("E24").Value = "  +0.89%  "

# Row 25 - Dai
$ws.Range("E25").Value = "  +0.08%  "

# Row 26 - PEPE
$ws.Range("E26").Value = "  +3.96%  "

# Row 27 - RenderToken
$ws.Range("E27").Value = "  +3.21%  "

# Row 28 - InternetComputer(DFINITY)
$ws.Range("D28").Value = "9.32"
$ws.Range("E28").Value = "  +2.19%  "

# Row 29 - PancakeSwap
$ws.Range("E29").Value = "  -0.42%  "

# Row 30 - Binance-PegBSC-USD
$ws.Range("E30").Value = "  +0.10%  "

# Row 31 - was RenzoRestakedETH, now Kaspa
$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D31").Value = "0.160"
$ws.Range("E31").Value = "  +3.84%  "

# Row 32 - was Kaspa, now RenzoRestakedETH
$ws.Range("B32").Value = "RenzoRestakedETH"
$ws.Range("C32").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D32").Value = "3.594.47"
$ws.Range("E32").Value = "  +1.08%  "

# Row 33 - Fetch.AI
$ws.Range("D33").Value = "1.47"
$ws.Range("E33").Value = "  -0.99%  "

# Row 34 - EthereumClassic
$ws.Range("D34").Value = "25.45"
$ws.Range("E34").Value = "  -0.32%  "

# Row 35 - Aptos
$ws.Range("D35").Value = "7.85"
$ws.Range("E35").Value = "  -0.78%  "

# Row 36 - USDe
$ws.Range("E36").Value = "  +0.00%  "

# Row 37 - NEARProtocol
$ws.Range("E37").Value = "  +0.27%  "

# Row 38 - ImmutableX
$ws.Range("E38").Value = "  -2.27%  "

# Row 39 - Monero
$ws.Range("D39").Value = "177.62"
$ws.Range("E39").Value = "  +4.42%  "

# Row 40 - Hedera
$ws.Range("E40").Value = "  +0.29%  "

# Row 41 - Filecoin
$ws.Range("D41").Value = "5.23"
$ws.Range("E41").Value = "  +0.44%  "

# Row 42 - Mantle
$ws.Range("E42").Value = "  +0.25%  "

# Row 43 - Stacks
$ws.Range("E43").Value = "  -1.53%  "

# Row 44 - dogwifhat
$ws.Range("E44").Value = "  +8.16%  "

# Row 45 - FirstDigitalUSD
$ws.Range("E45").Value = "  +0.04%  "

# Row 46 - ONDO
$ws.Range("E46").Value = "  -1.49%  "

# Row 47 - was InjectiveProtocol, now EnergySwap
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "24.35"
$ws.Range("E47").Value = "  +3.48%  "

# Row 48 - was EnergySwap, now InjectiveProtocol
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "25.01"
$ws.Range("E48").Value = "  -3.82%  "

# Row 49 - Cosmos
$ws.Range("D49").Value = "7.19"
$ws.Range("E49").Value = "  +0.62%  "

# Row 50 - SuiNetwork
$ws.Range("D50").Value = "0.952"
$ws.Range("E50").Value = "  +0.35%  "

# Row 51 - TheGraph
$ws.Range("E51").Value = "  -1.34%  "

# Restore the column's style so these cells keep the default (unstyled)
# appearance instead of retaining the temporary text number format.
$ws.Range("D2:D51").Style = "Normal"
